# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Prices in column D are plain text (not numbers) in the source sheet, so any
# value that Excel would otherwise auto-parse as a number is written through
# Set-TextValue, which forces the Text format, writes it, then pastes the
# (unstyled) format back from A1 so the cell's style index is left untouched,
# exactly matching the original (style-less) text cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $value
    $ws.Range("A1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)  # xlPasteFormats: restore default (unstyled) format
}

# Row 2: Bitcoin
$ws.Range("D2").Value = "89.560.15"
$ws.Range("E2").Value = "  -1.21%  "
# Row 3: Ethereum
$ws.Range("D3").Value = "3.069.51"
$ws.Range("E3").Value = "  -2.67%  "
# Row 4: TetherUSD
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.19%  "
# Row 5: Solana
Set-TextValue "D5" "234.58"
$ws.Range("E5").Value = "  +8.90%  "
# Row 6: BNB
Set-TextValue "D6" "617.94"
$ws.Range("E6").Value = "  -1.02%  "
# Row 7: XRP
Set-TextValue "D7" "1.07"
$ws.Range("E7").Value = "  -6.00%  "
# Row 8: Dogecoin
Set-TextValue "D8" "0.362"
$ws.Range("E8").Value = "  -1.10%  "
# Row 9: USDC
Set-TextValue "D9" "0.999"
$ws.Range("E9").Value = "  -0.08%  "
# Row 10: LidoStakedEther
$ws.Range("D10").Value = "3.072.56"
$ws.Range("E10").Value = "  -2.58%  "
# Row 11: Cardano
Set-TextValue "D11" "0.707"
$ws.Range("E11").Value = "  -5.82%  "
# Row 12: TRON
$ws.Range("E12").Value = "  -1.40%  "
# Row 13: ShibaInu
Set-TextValue "D13" "0.0000246"
$ws.Range("E13").Value = "  +0.66%  "
# Row 14: Avalanche
Set-TextValue "D14" "34.78"
$ws.Range("E14").Value = "  -1.36%  "
# Row 15: WrappedBTC
$ws.Range("D15").Value = "89.343.09"
$ws.Range("E15").Value = "  -1.29%  "
# Row 16: Toncoin
Set-TextValue "D16" "5.34"
$ws.Range("E16").Value = "  -6.11%  "
# Row 17: WrappedliquidstakedEther2.0
$ws.Range("D17").Value = "3.637.36"
$ws.Range("E17").Value = "  -2.69%  "
# Row 18: WrappedEther
$ws.Range("D18").Value = "3.074.86"
$ws.Range("E18").Value = "  -3.09%  "
# Row 19: SuiNetwork
Set-TextValue "D19" "3.73"
$ws.Range("E19").Value = "  -0.33%  "
# Row 20: PEPE
$ws.Range("E20").Value = "  +0.50%  "
# Row 21: Chainlink
Set-TextValue "D21" "13.68"
$ws.Range("E21").Value = "  -6.89%  "
# Row 22: BitcoinCash
Set-TextValue "D22" "429.11"
$ws.Range("E22").Value = "  -8.35%  "
# Row 23: Polkadot
Set-TextValue "D23" "5.37"
$ws.Range("E23").Value = "  +3.64%  "
# Row 24: Uniswap
Set-TextValue "D24" "8.65"
$ws.Range("E24").Value = "  -5.48%  "
# Row 25: NEARProtocol
Set-TextValue "D25" "5.53"
$ws.Range("E25").Value = "  -6.08%  "
# Row 26: Litecoin
Set-TextValue "D26" "87.12"
$ws.Range("E26").Value = "  -9.85%  "
# Row 27: Aptos
Set-TextValue "D27" "11.65"
$ws.Range("E27").Value = "  -5.47%  "
# Row 28: WrappedeETH
$ws.Range("D28").Value = "3.237.07"
# Row 29: Dai
$ws.Range("E29").Value = "  +0.03%  "
# Row 30: Binance-PegBSC-USD
Set-TextValue "D30" "1.07"
$ws.Range("E30").Value = "  +6.93%  "
# Row 31: InternetComputer(DFINITY)
Set-TextValue "D31" "8.99"
$ws.Range("E31").Value = "  -3.01%  "
# Row 32: Cronos
$ws.Range("E32").Value = "  -5.74%  "
# Row 33: Stellar
Set-TextValue "D33" "0.199"
$ws.Range("E33").Value = "  -9.88%  "
# Row 34: EthereumClassic
Set-TextValue "D34" "25.54"
$ws.Range("E34").Value = "  -5.54%  "
# Row 35: Kaspa
Set-TextValue "D35" "0.150"
$ws.Range("E35").Value = "  +3.29%  "
# Row 36: MantraDAO
Set-TextValue "D36" "4.02"
$ws.Range("E36").Value = "  +64.62%  "
# Row 37: RenderToken
Set-TextValue "D37" "6.98"
$ws.Range("E37").Value = "  -0.50%  "
# Row 38: Bittensor
Set-TextValue "D38" "489.26"
$ws.Range("E38").Value = "  -5.55%  "
# Row 39: dogwifhat
Set-TextValue "D39" "3.60"
$ws.Range("E39").Value = "  -0.25%  "
# Row 40: PancakeSwap
$ws.Range("E40").Value = "  -3.61%  "
# Row 41: Hedera (was Fetch.AI)
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D41" "0.0900"
$ws.Range("E41").Value = "  -2.32%  "
# Row 42: Fetch.AI (was Hedera)
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D42" "1.25"
$ws.Range("E42").Value = "  -7.00%  "
# Row 43: WhiteBITCoin
$ws.Range("E43").Value = "  -0.75%  "
# Row 44: USDe
$ws.Range("E44").Value = "  -0.02%  "
# Row 45: PolygonEcosystemToken
Set-TextValue "D45" "0.396"
$ws.Range("E45").Value = "  -7.52%  "
# Row 46: Monero
Set-TextValue "D46" "157.32"
$ws.Range("E46").Value = "  +4.41%  "
# Row 47: Stacks
Set-TextValue "D47" "1.84"
$ws.Range("E47").Value = "  -7.46%  "
# Row 48: ARBITRUM
Set-TextValue "D48" "0.671"
$ws.Range("E48").Value = "  -9.37%  "
# Row 49: OKB
Set-TextValue "D49" "44.15"
$ws.Range("E49").Value = "  -2.33%  "
# Row 50: FirstDigitalUSD
Set-TextValue "D50" "0.999"
$ws.Range("E50").Value = "  -0.38%  "
# Row 51: ImmutableX
$ws.Range("E51").Value = "  -6.59%  "

$excel.CutCopyMode = $false
